# Update the dSF (column F) values for the rows that were re-pulled / recalculated.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    11 = -1
    18 = 1
    20 = -1
    24 = -1
    26 = -10
    33 = 3
    35 = -2
    45 = -2
    53 = -1
    60 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
